$wb = $excel.ActiveWorkbook

# --- Sheet "Logs": append row 9 with the new test-mail entry ---
$ws = $wb.Worksheets.Item("Logs")

$ws.Range("A9").Value = "Is dit artikel nog op voorraad?"
$ws.Range("B9").Value = "mailmind.test@zohomail.eu"
$ws.Range("C9").Value = "Testmail #7: Is dit artikel nog op voorraad?"
$ws.Range("D9").Value = "Productinformatie"
$ws.Range("E9").Value = "Dank voor je bericht. We pakken dit intern op en houden je op de hoogte."
$ws.Range("F9").Value = "2025-07-31 21:36:52"
$ws.Range("G9").Value = "Ja"
$ws.Range("H9").Value = "Ja"
$ws.Range("I9").Value = "Nee"
$ws.Range("J9").Value = "Nee"

# --- Extend the conditional-formatting ranges from row 8 to row 9 ---
$colsToExtend = @("D", "G", "H", "I", "J")
foreach ($col in $colsToExtend) {
    $oldRange = $ws.Range($col + "2:" + $col + "8")
    $newRange = $ws.Range($col + "2:" + $col + "9")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# --- Sheet "Dashboard": update the category summary (re-sorted by count) ---
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A4").Value = "Productinformatie"
$dash.Range("B4").Value = 2
$dash.Range("A5").Value = "Bestelling / Levering"
$dash.Range("B5").Value = 1

$wb.Save()
